# Auto-generated PowerShell COM-interop script
# Applies "Add data for 2025-05-17" edits to violent-crime-full-year.xlsx
$wb = $excel.ActiveWorkbook

# --- Citywide Totals ---
$ws = $wb.Worksheets.Item("Citywide Totals")
$ws.Range("L2").Value = 2307
$ws.Range("L3").Value = 2350
$ws.Range("F4").Value = 1928
$ws.Range("H4").Value = 1755
$ws.Range("I4").Value = 1837
$ws.Range("K4").Value = 1765
$ws.Range("L4").Value = 636
$ws.Range("L6").Value = 2115
$ws.Range("F7").Value = 24121
$ws.Range("H7").Value = 26068
$ws.Range("I7").Value = 26305
$ws.Range("K7").Value = 27557
$ws.Range("L7").Value = 7548

# --- By Neighborhood ---
$ws = $wb.Worksheets.Item("By Neighborhood")
$ws.Range("L4").Value = 32
$ws.Range("L5").Value = 25
$ws.Range("L6").Value = 59
$ws.Range("L7").Value = 246
$ws.Range("L8").Value = 479
$ws.Range("L9").Value = 50
$ws.Range("L11").Value = 133
$ws.Range("L15").Value = 53
$ws.Range("L19").Value = 217
$ws.Range("L20").Value = 191
$ws.Range("L23").Value = 78
$ws.Range("L28").Value = 4
$ws.Range("L29").Value = 380
$ws.Range("L33").Value = 346
$ws.Range("L36").Value = 108
$ws.Range("L37").Value = 275
$ws.Range("L43").Value = 63
$ws.Range("L44").Value = 55
$ws.Range("L47").Value = 58
$ws.Range("L48").Value = 105
$ws.Range("L49").Value = 43
$ws.Range("L51").Value = 86
$ws.Range("L52").Value = 149
$ws.Range("L54").Value = 155
$ws.Range("L55").Value = 69
$ws.Range("L60").Value = 42
$ws.Range("F63").Value = 213
$ws.Range("H63").Value = 304
$ws.Range("I63").Value = 260
$ws.Range("K63").Value = 156
$ws.Range("L63").Value = 24
$ws.Range("L65").Value = 143
$ws.Range("L67").Value = 277
$ws.Range("L73").Value = 60
$ws.Range("L76").Value = 84
$ws.Range("L77").Value = 47
$ws.Range("L79").Value = 205
$ws.Range("L81").Value = 7
$ws.Range("L83").Value = 184
$ws.Range("L84").Value = 74
$ws.Range("L85").Value = 393
$ws.Range("L88").Value = 102
$ws.Range("L90").Value = 76
$ws.Range("L93").Value = 40
$ws.Range("L96").Value = 73
$ws.Range("L98").Value = 54
$ws.Range("L100").Value = 11
$ws.Range("F101").Value = 24121
$ws.Range("H101").Value = 26068
$ws.Range("I101").Value = 26305
$ws.Range("K101").Value = 27557
$ws.Range("L101").Value = 7548

# --- West Ridge ---
$ws = $wb.Worksheets.Item("West Ridge")
$ws.Range("L4").Value = 9
$ws.Range("L7").Value = 73

# --- Auburn Gresham ---
$ws = $wb.Worksheets.Item("Auburn Gresham")
$ws.Range("L3").Value = 80
$ws.Range("L6").Value = 69
$ws.Range("L7").Value = 246

# --- Belmont Cragin ---
$ws = $wb.Worksheets.Item("Belmont Cragin")
$ws.Range("L2").Value = 46
$ws.Range("L7").Value = 133

# --- South Shore ---
$ws = $wb.Worksheets.Item("South Shore")
$ws.Range("L3").Value = 162
$ws.Range("L4").Value = 31
$ws.Range("L7").Value = 393

# --- Little Village ---
$ws = $wb.Worksheets.Item("Little Village")
$ws.Range("L2").Value = 46
$ws.Range("L6").Value = 43
$ws.Range("L7").Value = 149

# --- Austin ---
$ws = $wb.Worksheets.Item("Austin")
$ws.Range("L3").Value = 163
$ws.Range("L6").Value = 128
$ws.Range("L7").Value = 479

# --- South Chicago ---
$ws = $wb.Worksheets.Item("South Chicago")
$ws.Range("L3").Value = 75
$ws.Range("L6").Value = 42
$ws.Range("L7").Value = 184

# --- Garfield Park ---
$ws = $wb.Worksheets.Item("Garfield Park")
$ws.Range("L2").Value = 92
$ws.Range("L3").Value = 110
$ws.Range("L6").Value = 119
$ws.Range("L7").Value = 346

# --- Grand Crossing ---
$ws = $wb.Worksheets.Item("Grand Crossing")
$ws.Range("L2").Value = 82
$ws.Range("L3").Value = 82
$ws.Range("L7").Value = 275

# --- New City ---
$ws = $wb.Worksheets.Item("New City")
$ws.Range("L2").Value = 53
$ws.Range("L7").Value = 143

# --- North Lawndale ---
$ws = $wb.Worksheets.Item("North Lawndale")
$ws.Range("L2").Value = 81
$ws.Range("L3").Value = 94
$ws.Range("L7").Value = 277

# --- South Deering ---
$ws = $wb.Worksheets.Item("South Deering")
$ws.Range("L6").Value = 17
$ws.Range("L7").Value = 74

# --- Lincoln Park ---
$ws = $wb.Worksheets.Item("Lincoln Park")
$ws.Range("L2").Value = 14
$ws.Range("L7").Value = 43

# --- Loop ---
$ws = $wb.Worksheets.Item("Loop")
$ws.Range("L3").Value = 29
$ws.Range("L6").Value = 81
$ws.Range("L7").Value = 155

# --- Englewood ---
$ws = $wb.Worksheets.Item("Englewood")
$ws.Range("L2").Value = 124
$ws.Range("L3").Value = 136
$ws.Range("L7").Value = 380

# --- Lake View ---
$ws = $wb.Worksheets.Item("Lake View")
$ws.Range("L6").Value = 42
$ws.Range("L7").Value = 105

# --- Chatham ---
$ws = $wb.Worksheets.Item("Chatham")
$ws.Range("L2").Value = 71
$ws.Range("L3").Value = 67
$ws.Range("L6").Value = 69
$ws.Range("L7").Value = 217

# --- Irving Park ---
$ws = $wb.Worksheets.Item("Irving Park")
$ws.Range("L2").Value = 22
$ws.Range("L3").Value = 15
$ws.Range("L7").Value = 55

# --- River North ---
$ws = $wb.Worksheets.Item("River North")
$ws.Range("L6").Value = 40
$ws.Range("L7").Value = 84

# --- Ashburn ---
$ws = $wb.Worksheets.Item("Ashburn")
$ws.Range("L6").Value = 12
$ws.Range("L7").Value = 59

# --- Lower West Side ---
$ws = $wb.Worksheets.Item("Lower West Side")
$ws.Range("L3").Value = 25
$ws.Range("L7").Value = 69

# --- Douglas ---
$ws = $wb.Worksheets.Item("Douglas")
$ws.Range("L6").Value = 16
$ws.Range("L7").Value = 78

# --- Roseland ---
$ws = $wb.Worksheets.Item("Roseland")
$ws.Range("L2").Value = 66
$ws.Range("L7").Value = 205

# --- Chicago Lawn ---
$ws = $wb.Worksheets.Item("Chicago Lawn")
$ws.Range("L2").Value = 61
$ws.Range("L7").Value = 191

# --- Grand Boulevard ---
$ws = $wb.Worksheets.Item("Grand Boulevard")
$ws.Range("L2").Value = 46
$ws.Range("L3").Value = 26
$ws.Range("L4").Value = 9
$ws.Range("L7").Value = 108

# --- West Lawn ---
$ws = $wb.Worksheets.Item("West Lawn")
$ws.Range("L2").Value = 15
$ws.Range("L7").Value = 40

# --- Wrigleyville ---
$ws = $wb.Worksheets.Item("Wrigleyville")
$ws.Range("L3").Value = 5
$ws.Range("L7").Value = 11

# --- Kenwood ---
$ws = $wb.Worksheets.Item("Kenwood")
$ws.Range("L3").Value = 20
$ws.Range("L7").Value = 58

# --- Brighton Park ---
$ws = $wb.Worksheets.Item("Brighton Park")
$ws.Range("L2").Value = 20
$ws.Range("L7").Value = 53

# --- Wicker Park ---
$ws = $wb.Worksheets.Item("Wicker Park")
$ws.Range("L6").Value = 29
$ws.Range("L7").Value = 54

# --- Avalon Park ---
$ws = $wb.Worksheets.Item("Avalon Park")
$ws.Range("L3").Value = 22
$ws.Range("L7").Value = 50

# --- Portage Park ---
$ws = $wb.Worksheets.Item("Portage Park")
$ws.Range("L2").Value = 25
$ws.Range("L7").Value = 60

# --- United Center ---
$ws = $wb.Worksheets.Item("United Center")
$ws.Range("L2").Value = 27
$ws.Range("L7").Value = 102

# --- Armour Square ---
$ws = $wb.Worksheets.Item("Armour Square")
$ws.Range("L2").Value = 7
$ws.Range("L7").Value = 25

# --- Washington Heights ---
$ws = $wb.Worksheets.Item("Washington Heights")
$ws.Range("L3").Value = 22
$ws.Range("L4").Value = 6
$ws.Range("L7").Value = 76

# --- Little Italy, UIC ---
$ws = $wb.Worksheets.Item("Little Italy, UIC")
$ws.Range("L2").Value = 24
$ws.Range("L4").Value = 10
$ws.Range("L7").Value = 86

# --- Morgan Park ---
$ws = $wb.Worksheets.Item("Morgan Park")
$ws.Range("L3").Value = 14
$ws.Range("L7").Value = 42

# --- Hyde Park ---
$ws = $wb.Worksheets.Item("Hyde Park")
$ws.Range("L3").Value = 17
$ws.Range("L4").Value = 14
$ws.Range("L6").Value = 20
$ws.Range("L7").Value = 63

# --- Riverdale ---
$ws = $wb.Worksheets.Item("Riverdale")
$ws.Range("L3").Value = 18
$ws.Range("L7").Value = 47

# --- Archer Heights ---
$ws = $wb.Worksheets.Item("Archer Heights")
$ws.Range("L2").Value = 11
$ws.Range("L3").Value = 10
$ws.Range("L7").Value = 32

# --- Sauganash,Forest Glen ---
$ws = $wb.Worksheets.Item("Sauganash,Forest Glen")
$ws.Range("L3").Value = 2
$ws.Range("L7").Value = 7

# --- Edison Park ---
$ws = $wb.Worksheets.Item("Edison Park")
$ws.Range("L6").ClearContents()
$ws.Range("L7").Value = 4
